$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# 1) Column B (username) for rows 1-42 changes from "aline" to "vanessa"
$ws.Range("B1:B42").Value = "vanessa"

# 2) Row 43 keeps the same grant text, but the username changes from "aline" to "evaldo"
#    (this row now belongs to the new "evaldo" grant block)
$ws.Cells.Item(43, 2).Value = "evaldo"

# 3) Copy formatting of row 43 down into the newly-added rows 44:49
$ws.Range("A43:D43").Copy() | Out-Null
$ws.Range("A44:D49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4) Fill in the new grant rows for user "evaldo"
$ws.Cells.Item(44, 1).Value = 'GRANT ALL ON SEQUENCE tratamento.sq_risco_rnado_pcnt TO'
$ws.Cells.Item(44, 2).Value = 'evaldo'
$ws.Cells.Item(44, 3).Value = ';'
$ws.Cells.Item(44, 4).Formula = "=A44&`" `"&B44&`" `"&C44"

$ws.Cells.Item(45, 1).Value = 'GRANT ALL ON SEQUENCE tratamento.sq_hstr_pnel_mapa_risco TO'
$ws.Cells.Item(45, 2).Value = 'evaldo'
$ws.Cells.Item(45, 3).Value = ';'
$ws.Cells.Item(45, 4).Formula = "=A45&`" `"&B45&`" `"&C45"

$ws.Cells.Item(46, 1).Value = 'GRANT ALL ON SEQUENCE tratamento.sq_hstr_obs_pnel_mapa_risco TO'
$ws.Cells.Item(46, 2).Value = 'evaldo'
$ws.Cells.Item(46, 3).Value = ';'
$ws.Cells.Item(46, 4).Formula = "=A46&`" `"&B46&`" `"&C46"

$ws.Cells.Item(47, 1).Value = 'GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_hstr_obs_pnel_mapa_risco TO'
$ws.Cells.Item(47, 2).Value = 'evaldo'
$ws.Cells.Item(47, 3).Value = ';'
$ws.Cells.Item(47, 4).Formula = "=A47&`" `"&B47&`" `"&C47"

$ws.Cells.Item(48, 1).Value = 'GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_hstr_pnel_mapa_risco TO'
$ws.Cells.Item(48, 2).Value = 'evaldo'
$ws.Cells.Item(48, 3).Value = ';'
$ws.Cells.Item(48, 4).Formula = "=A48&`" `"&B48&`" `"&C48"

$ws.Cells.Item(49, 1).Value = 'GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_risco_rnado_pcnt TO'
$ws.Cells.Item(49, 2).Value = 'evaldo'
$ws.Cells.Item(49, 3).Value = ';'
$ws.Cells.Item(49, 4).Formula = "=A49&`" `"&B49&`" `"&C49"

